$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (pushes old row 11 "Daily Total" down to row 12)
$ws.Rows("11:11").Insert()

# New row 11: "Charter Work" task with 0.5 hours on Tuesday (column C) and 0.5 weekly total (column I)
$ws.Range("A11").Value = "Charter Work"
$ws.Range("C11").Value = 0.5
$ws.Range("I11").Value = 0.5

# Update the Daily Total row (now row 12) Tuesday column to include the new 0.5 hours
$ws.Range("C12").Value = 0.5

# Update selection to match target state
$ws.Range("I12").Select()
